$d = $word.ActiveDocument

# 1 & 4) "January 27, 2022" -> "January 28, 2022" (appears twice: arraignment date and
#    the "fines and costs shall be paid in full by" date). Both occurrences change.
$d.Content.Find.Execute("January 27, 2022", $true, $false, $false, $false, $false, `
    $true, 1, $false, "January 28, 2022", 2)

# 2 & 3) In the charges table, the "Plea" row changes from "Guilty"/"Guilty" to
#    "No Contest"/"No Contest". The "Finding" row below it (also "Guilty"/"Guilty")
#    must stay untouched, so target the Plea row's two cells directly.
$table = $d.Tables.Item(1)
$pleaRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $label = $table.Cell($i, 1).Range.Text
    if ($label -like "Plea*") {
        $pleaRow = $i
        break
    }
}
$cell2 = $table.Cell($pleaRow, 2)
$cell2.Range.Find.Execute("Guilty", $true, $false, $false, $false, $false, `
    $true, 0, $false, "No Contest", 1)
$cell3 = $table.Cell($pleaRow, 3)
$cell3.Range.Find.Execute("Guilty", $true, $false, $false, $false, $false, `
    $true, 0, $false, "No Contest", 1)

# 5) Proof of Financial Responsibility paragraph: the sentence about the Defendant
#    showing proof is replaced with a sentence saying the Defendant did NOT show
#    proof (and may still show it to the Clerk before the matter goes to the BMV).
#    The sentence is preceded by a single space run that collapses away once its
#    text is cleared, so remove that leading space and then swap in the new text.
$old = "The Defendant showed proof of financial responsibility at the time of the offense."
$new = "The Defendant did not show proof of financial responsibility at the time of the offense or during the proceeding, but may show proof to Clerk of Court at any time prior to the submission of this matter to the Ohio Bureau of Motor Vehicles."

$findRng = $d.Content
$findRng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $findRng.Start
$end = $findRng.End

$spaceRng = $d.Range($start - 1, $start)
$spaceRng.Text = ""

$textRng = $d.Range($start - 1, $end - 1)
$textRng.Text = $new
